# Refresh the cryptocurrency Price / Volume(1h) columns with newly scraped values.
# Price cells are forced to plain text so values like "65.046.88" or "1.00"
# are preserved exactly instead of being reinterpreted as numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.046.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.101.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.101.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("E11").Value = "  +4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.614.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.205.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.109.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  +6.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.16%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "470.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0412"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +20.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.001.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.90%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.112"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0531"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.92%  "

Write-Host "Updated cryptos list"
